$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Add a new (unbound) column to the right of the query table, matching
# Excel's "insert column" behaviour next to a query-table-backed table.
$newCol = $tbl.ListColumns.Add()

# Name the new column via its header cell (ListColumn.Name assignment
# doesn't propagate reliably, but writing the header cell value does).
$headerCell = $newCol.Range.Cells.Item(1, 1)
$headerCell.Value = "QTY per 100"

# Populate "QTY per 100" = "Quantity Per PCB" * 100 for every data row.
$dataRange = $newCol.DataBodyRange
for ($r = 1; $r -le $dataRange.Rows.Count; $r++) {
    $qtyPerPcb = $ws.Cells.Item($r + 1, 8).Value()
    $dataRange.Cells.Item($r, 1).Value = $qtyPerPcb * 100
}

# Match the new column's visual style (s="1") to the rest of the table body.
$dataRange.Style = $ws.Range("H2").Style

# New column width, like the rest of the BOM columns.
$ws.Columns.Item(9).ColumnWidth = 11.85546875

# Reflect the freshly-added column in the view/selection, as Excel would
# after inserting + selecting the new column's data body.
$ws.Application.ActiveWindow.ScrollColumn = 2
$dataRange.Select()
